$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Stack trace line-number updates (POI 3.17.0 -> 4.0.1 / tooling upgrade)
Replace-Text "JavaMethodService.internalInvoke(JavaMethodService.java:163)" "JavaMethodService.internalInvoke(JavaMethodService.java:162)"
Replace-Text "AbstractService.invoke(AbstractService.java:136)" "AbstractService.invoke(AbstractService.java:135)"
Replace-Text "EvaluationServices.call(EvaluationServices.java:168)" "EvaluationServices.call(EvaluationServices.java:172)"
Replace-Text "AstSwitch.doSwitch(AstSwitch.java:118)" "AstSwitch.doSwitch(AstSwitch.java:119)"
Replace-Text "M2DocEvaluator.caseLink(M2DocEvaluator.java:1703)" "M2DocEvaluator.caseLink(M2DocEvaluator.java:1705)"
Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)"
Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)"
Replace-Text "GeneratedMethodAccessor74" "GeneratedMethodAccessor75"

# Replace the tail of the stack trace (Maven/Tycho/Equinox launcher frames)
# with the Eclipse JDT JUnit runner frames.
$oldTail = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n" +
    "	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n" +
    "	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
    "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
    "	at java.lang.reflect.Method.invoke(Method.java:498)`n" +
    "	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n" +
    "	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n" +
    "	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n" +
    "	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n" +
    "	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
    "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
    "	at java.lang.reflect.Method.invoke(Method.java:498)`n" +
    "	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n" +
    "	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n" +
    "	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n" +
    "	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n" +
    "	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n" +
    "	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
    "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
    "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
    "	at java.lang.reflect.Method.invoke(Method.java:498)`n" +
    "	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n" +
    "	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n" +
    "	at org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n" +
    "	at org.eclipse.equinox.launcher.Main.main(Main.java:1471)`n"

$newTail = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" +
    "	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" +
    "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" +
    "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" +
    "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" +
    "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"

Replace-Text $oldTail $newTail

# NOTE: the diff also shows the bookmark's internal w:id (a generated
# "random" 39-digit token) and the w:rsidR token stamped on the REF
# field's runs changing value. Both are purely-internal OOXML
# bookkeeping attributes (not part of any run's visible text), and are
# not exposed as writable properties anywhere on the Word object model
# (Bookmark has no settable Id, and there is no RsidR property), so
# they cannot be targeted with Find/Replace or any other COM call --
# Word (and this interop layer) manage/regenerate them on save.

Write-Host "done"
